# ydf_treeple_scaling.xlsx — "Update MIGHT notebook to make it identical to YDF"
#
# Fills in the previously-missing Treeple timing numbers (n_attributes = 160
# table, rows for n=2000/4000/8000) so the lower ("Treeple") table has the
# same shape as the upper ("YDF") table, and leaves the selection on the
# last cell touched (I17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15 (n=2000) ---------------------------------------------------
# The row already held 5 values in B15:F15 (for d=320..2048) but was
# missing the d=160 column and the d=4096 column. Shift the existing
# values one column to the right (B->C, C->D, D->E, E->F, F->G) to make
# room for the new first data point, then fill in the new first and last
# cells.
$ws.Range("G15").Value = $ws.Range("F15").Value()
$ws.Range("F15").Value = $ws.Range("E15").Value()
$ws.Range("E15").Value = $ws.Range("D15").Value()
$ws.Range("D15").Value = $ws.Range("C15").Value()
$ws.Range("C15").Value = $ws.Range("B15").Value()
$ws.Range("B15").Value = 54.9833
$ws.Range("H15").Value = 802.4954

# --- Row 16 (n=4000) — previously empty except n value -----------------
$ws.Range("B16").Value = 140.0341
$ws.Range("C16").Value = 166.26
$ws.Range("D16").Value = 222.1643
$ws.Range("E16").Value = 296.4236
$ws.Range("F16").Value = 480.165
$ws.Range("G16").Value = 859.6107
$ws.Range("H16").Value = 1610.7064

# --- Row 17 (n=8000) — previously empty except n value ------------------
$ws.Range("B17").Value = 341.1712
$ws.Range("C17").Value = 389.9193
$ws.Range("D17").Value = 496.8179
$ws.Range("E17").Value = 643.755
$ws.Range("F17").Value = 1019.6019
$ws.Range("G17").Value = 1776.1397
$ws.Range("H17").Value = 3272.4283

# --- Match the author's final selection ---------------------------------
$ws.Range("I17").Select()
